$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row at row 3, pushing the existing "CHORA" / "RICO" rows down by one.
$ws.Rows.Item(3).Insert()

# Fill in the data for the newly inserted student record.
$ws.Range("A3").Value = 19330051920235
$ws.Range("B3").Value = "HERAS"
$ws.Range("C3").Value = "LOPEZ"
$ws.Range("D3").Value = "CESAR ENRIQUE"
$ws.Range("E3").Value = "CONSTRUYE BASES DE DATOS PARA APLICACIONES WEB"
$ws.Range("F3").Value = "5APM"
$ws.Range("G3").Value = 6
